# Add data for 2022-10-31: updates the "October" (most-recent-month) column
# for each neighborhood row, plus the header date label and sheet name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the report "as of" date in the sheet name and the column header text.
$ws.Name = "Through 2022-10-23"
$ws.Range("B1").Value = "October 2022 (through October 23)"

# Per-neighborhood October counts, keyed by cell reference.
$updates = @{
    "V2"   = 16
    "AF2"  = 5
    "B3"   = 5
    "L3"   = 7
    "V3"   = 2
    "AF3"  = 4
    "BT4"  = 3
    "L5"   = 13
    "L6"   = 12
    "L8"   = 4
    "AF8"  = 3
    "L9"   = 5
    "B12"  = 5
    "V13"  = 1
    "AZ13" = 2
    "L15"  = 1
    "B18"  = 4
    "BJ19" = 1
    "BJ22" = 1
    "B23"  = 2
    "V23"  = 2
    "AZ24" = 2
    "L25"  = 2
    "B27"  = 1
    "B29"  = 1
    "V36"  = 1
    "B37"  = 1
    "L38"  = 1
    "V46"  = 4
    "AZ47" = 1
    "V49"  = 1
    "B58"  = 1
    "L69"  = 1
    "AP69" = 1
    "V97"  = 1
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
